$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.635.70"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "3.626.07"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "204.40"
$ws.Range("E5").Value = "  +10.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "569.49"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("D7").Value = "3.619.22"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.625"
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.680"
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "61.88"
$ws.Range("E11").Value = "  +17.50%  "
$ws.Range("E12").Value = "  +4.95%  "
$ws.Range("E13").Value = "  +11.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.10"
$ws.Range("E14").Value = "  +3.47%  "
$ws.Range("D15").Value = "4.197.09"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").Value = "3.617.22"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.09"
$ws.Range("E18").Value = "  +4.56%  "
$ws.Range("D19").Value = "68.256.68"
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.45"
$ws.Range("E20").Value = "  +2.54%  "
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "405.76"
$ws.Range("E22").Value = "  +3.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.38"
$ws.Range("E23").Value = "  +20.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.20"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.87"
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.96"
$ws.Range("E26").Value = "  +2.26%  "
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.95"
$ws.Range("E28").Value = "  +12.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.13"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.39"
$ws.Range("E30").Value = "  +5.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.07"
$ws.Range("E31").Value = "  +14.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.75"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "668.85"
$ws.Range("E33").Value = "  +6.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.36"
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.72"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.42"
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.422"
$ws.Range("E38").Value = "  +6.50%  "
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").Value = "0.0₃0789"
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("E41").Value = "  +14.47%  "
$ws.Range("D42").Value = "3.220.49"
$ws.Range("E42").Value = "  +8.47%  "
$ws.Range("E43").Value = "  +4.08%  "
$ws.Range("E44").Value = "  +11.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.996"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.97"
$ws.Range("E46").Value = "  +25.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.89"
$ws.Range("E47").Value = "  +15.72%  "
$ws.Range("E48").Value = "  +3.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.92"
$ws.Range("E49").Value = "  +5.58%  "
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.09"
$ws.Range("E51").Value = "  -0.13%  "
